# UX Upgrade: Hero Search, SEO, and 20+ Categorized Posts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had 4 columns (topic, category, banner_type, prompt) across
# 11 rows. The new layout drops banner_type, reorders the remaining columns
# to (topic, prompt, category), and grows to 20 data rows (21 rows incl.
# header). Drop column D completely (content + formatting), and clear out
# the old data rows before repopulating A:C with the new content. The
# header row's existing formatting (bold/border/centered) is left as-is.
$ws.Range("D1:D11").Clear()
$ws.Range("A2:C11").ClearContents()

$ws.Cells.Item(1, 1).Value = "topic"
$ws.Cells.Item(1, 2).Value = "prompt"
$ws.Cells.Item(1, 3).Value = "category"

$topics = @(
    "Instagram Viral Caption Generator",
    "Python Code Debugger",
    "SEO Blog Post Writer",
    "Email Cold Outreach Template",
    "Midjourney Portrait Prompt",
    "YouTube Video Script Hook",
    "React Component Generator",
    "Business Plan Executive Summary",
    "Fitness Workout Plan AI",
    "Meditation Script Generator",
    "TikTok Trend Ideas",
    "SQL Query Optimizer",
    "Resume Bullet Point Polisher",
    "Recipe Generator from Ingredients",
    "Startup Pitch Deck Outliner",
    "Twitter Thread Creator",
    "Excel Formula Explainer",
    "Novel Character Backstory",
    "Product Description Improver",
    "Daily Journaling Prompts"
)

$categories = @(
    "Marketing",
    "Coding",
    "Writing",
    "Business",
    "Art",
    "Marketing",
    "Coding",
    "Business",
    "Lifestyle",
    "Lifestyle",
    "Marketing",
    "Coding",
    "Business",
    "Lifestyle",
    "Business",
    "Marketing",
    "Coding",
    "Writing",
    "Marketing",
    "Lifestyle"
)

for ($i = 0; $i -lt $topics.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $topics[$i]
    $ws.Cells.Item($row, 2).Value = "Write a viral caption for..."
    $ws.Cells.Item($row, 3).Value = $categories[$i]
}

Write-Host "Applied 20 categorized rows"
